$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$app = $excel

# New device rows to append (group "32"), mirroring the pattern of the
# previous groups (rows 2-156) already present in the sheet.
$rows = @(
    @{ Row = 157; Id = 3000176; Name = "Finger Print Scanner 32"; Mac = "80-75-40-E8-CA-24"; Serial = "BS563Q2230824"; Dspec = 165 },
    @{ Row = 158; Id = 3000177; Name = "IRIS Scanner 32";          Mac = "0E-1A-14-4A-6D-3A"; Serial = "BS563Q2230825"; Dspec = 327 },
    @{ Row = 159; Id = 3000178; Name = "Web Camera 32";            Mac = "65-13-7F-0F-F7-53"; Serial = "BS563Q2230826"; Dspec = 736 },
    @{ Row = 160; Id = 3000179; Name = "Document Scanner 32";      Mac = "73-C4-DE-8E-C9-8D"; Serial = "BS563Q2230827"; Dspec = 801 },
    @{ Row = 161; Id = 3000180; Name = "Printer 32";                Mac = "EC-74-AB-E0-0F-38"; Serial = "BS563Q2230828"; Dspec = 920 }
)

# Column A (ids) first, matching how the new strings were introduced
# column-by-column (B, then C, then D) in the source workbook.
foreach ($r in $rows) {
    $ws.Cells.Item($r.Row, 1).Value = $r.Id
}
foreach ($r in $rows) {
    $ws.Cells.Item($r.Row, 2).Value = $r.Name
}
foreach ($r in $rows) {
    $ws.Cells.Item($r.Row, 3).Value = $r.Mac
}
foreach ($r in $rows) {
    $ws.Cells.Item($r.Row, 4).Value = $r.Serial
}
foreach ($r in $rows) {
    $i = $r.Row
    $ws.Cells.Item($i, 6).Value = $r.Dspec
    $ws.Cells.Item($i, 7).Value = "eng"
    $ws.Cells.Item($i, 8).Value = $true
    $ws.Cells.Item($i, 8).HorizontalAlignment = -4131
    $ws.Cells.Item($i, 9).Value = "superadmin"
    $ws.Cells.Item($i, 10).Value = "now()"
}

# Restore view state: scroll/selection ends up on columns K:XFD with the
# active cell around row 113, matching where the author left the sheet.
$ws.Activate() | Out-Null
$ws.Range("K1:XFD1048576").Select() | Out-Null
